# Add a new "2022" data column (P) to the right of the existing "2021"
# column (O), mirroring column O's per-cell formatting and filling in the
# 2022 figures for each region.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring column P's formatting in line with column O (same borders, fonts,
# number formats, etc. row by row) by copying each O cell onto its P
# neighbour; the values get overwritten right after. Only rows 4-14 hold
# table data in column O, so that's the only range that needs new P cells.
for ($r = 4; $r -le 14; $r++) {
    $ws.Range("O$r").Copy($ws.Range("P$r")) | Out-Null
}

# Header year for the new column
$ws.Range("P4").Value = 2022

# 2022 figures, one per region (rows 6-14), plus the summary row (5)
$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

# Leave the selection where the author's edit left it
$ws.Range("O21:O22").Select() | Out-Null
